# Update check_script/requirements.xlsx (Sheet1) to add new device/memory-tuning
# related requirement rows (lmk, dirty_ratio/dirty_background_ratio, readahead,
# io_scheduler) plus a new "device" column, per commit:
# "update some functions include lmk, readahead, writeback, ioscheduler compare"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 5 fresh rows right before the old "CPU_kernel_count" row (row 6),
# pushing the CPU/Temperature block down to rows 11-15. Excel inherits the
# formatting of the row above (row 5, style index 1) for the newly inserted
# rows, so no extra cell styles are introduced.
$ws.Rows("6:10").Insert()

# Populate the newly inserted rows with the new requirement checks.
$ws.Range("A6").Value = "lmk_minfree_levels"
$ws.Range("B6").Value = 950

$ws.Range("A7").Value = "dirty_ratio"
$ws.Range("B7").Value = 30

$ws.Range("A8").Value = "dirty_background_ratio"
$ws.Range("B8").Value = 10

$ws.Range("A9").Value = "readahead"
$ws.Range("B9").Value = 1024
$ws.Range("E9").Value = "sda"

$ws.Range("A10").Value = "io_scheduler"
$ws.Range("B10").Value = "cfp"
$ws.Range("E10").Value = "sda"

# New "device" header column.
$ws.Range("E1").Value = "device"

# The old rows 6-10 (now 11-15, e.g. CPU_kernel_count.. Temperature) used a
# cell style that duplicated the "horizontal+vertical center" style already
# used by column B. Re-apply that existing (shared) format via a
# format-only paste so we reuse the existing style entry instead of
# creating a brand-new one.
[void]$ws.Range("B1").Copy()
[void]$ws.Range("A11:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A grew wider to fit the new, longer requirement names.
$ws.Columns("A").ColumnWidth = 24.14

# Leave the cursor/selection on the last newly entered cell, like the
# original author would have after typing in the new data.
[void]$ws.Range("B10").Select()
